# "Author has been Added" - add two new locator rows (text_search / btn_search)
# to the ObjectR sheet, bold the new "Locator_Name" entries, and update the
# active-sheet/selection bookkeeping accordingly.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ObjectR
$ws2 = $wb.Worksheets.Item(2)   # TestD

# --- New data rows on ObjectR --------------------------------------------
# Row 13: text_search / XPATH / //input[@title='Search']
# Set column C before column A so the shared-string table is populated in
# the same order as the target workbook.
$ws1.Range("C13").Value = "//input[@title='Search']"
$ws1.Range("A13").Value = "text_search"
$ws1.Range("B13").Value = "XPATH"
$ws1.Range("A13").Font.Bold = $true

# Row 14: btn_search / XPATH / //div[@class='FPdoLc VlcLAe']//input[@value='Google Search']
$ws1.Range("A14").Value = "btn_search"
$ws1.Range("C14").Value = "//div[@class='FPdoLc VlcLAe']//input[@value='Google Search']"
$ws1.Range("B14").Value = "XPATH"
$ws1.Range("A14").Font.Bold = $true

# --- Active sheet / selection bookkeeping --------------------------------
# TestD loses tab-selection, its lingering selection moves to D17.
$ws2.Activate()
$ws2.Range("D17").Select()

# ObjectR becomes (and stays) the selected/active tab, selection moves to C17.
$ws1.Activate()
$ws1.Range("C17").Select()
